# "more layout done to PCB" - update sensor overview sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# wire counts (col C) updated for a few sensors; D (wire/p unit * units) is a
# shared formula and recalculates automatically.
$ws.Range("C5").Value = 3    # CAN: 2 -> 3 wires/unit
$ws.Range("C6").Value = 3    # RS232: 2 -> 3 wires/unit
$ws.Range("B17").Value = 2   # Brake Pressure: 1 -> 2 units
$ws.Range("C21").Value = 7   # H-bridge module: 8 -> 7 wires/unit

# New "AUX" sensor row (row 25) filled in with units/wires and the "needs Int" (H) count
$ws.Range("A25").Value = "AUX"
$ws.Range("B25").Value = 4
$ws.Range("C25").Value = 2
$ws.Range("H25").Value = 4

# Mark a few more sensors as needing an ADC channel (col J, shared string "x")
$ws.Range("J13").Value = "x"
$ws.Range("J16").Value = "x"
$ws.Range("J19").Value = "x"

# "Con pins left" formula budget changed from 78 to 68
$ws.Range("I35").Formula = "=68-I34"

# Minor column width tweaks for columns A and B
$ws.Columns.Item(1).ColumnWidth = 28.5
$ws.Columns.Item(2).ColumnWidth = 8.833333333333332

# Move the active selection to C23
$excel.Goto($ws.Range("C23"))
